# Add a new column BL with header "04-sep" and the corresponding data values,
# mirroring the formatting of the adjacent BK column (header style = BK1's style).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell BL1: same style as BK1 (text-formatted header) + the new date label
$ws.Range("BL1").Value = "04-sep"
$ws.Range("BL1").NumberFormat = $ws.Range("BK1").NumberFormat

# Data cells BL2:BL18
$values = @{
    2  = 0
    3  = 16.713348746320555
    4  = 13.968397319158779
    5  = 10.925234804063917
    6  = 0
    7  = 14.719281109091291
    8  = 10.207087334673719
    9  = 9.0560117471410742
    10 = 17.498656503306165
    11 = 10.764454819496196
    12 = 0
    13 = 10.115268498261946
    14 = 0
    15 = 0
    16 = 17.197013873821557
    17 = 0
    18 = 0
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 64).Value = $values[$row]
}

# Update the active selection, matching the post-edit state captured by Excel
$ws.Range("BN5").Select()
